$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on column D (price) cells before assigning values,
# since these values must remain literal text (e.g. "1.000", "240.90")
# and not be auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.851.80'
$ws.Range("E2").Value = '  -0.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.737.62'
$ws.Range("E3").Value = '  -0.09%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.90'
$ws.Range("E5").Value = '  +4.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5229'
$ws.Range("E7").Value = '  -0.42%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2745'
$ws.Range("E8").Value = '  -0.80%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06166'
$ws.Range("E9").Value = '  +0.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.741.38'
$ws.Range("E10").Value = '  +0.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07192'
$ws.Range("E11").Value = '  +1.44%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.94'
$ws.Range("E12").Value = '  -1.25%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6426'
$ws.Range("E13").Value = '  -0.24%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.622'

$ws.Range("E15").Value = '  +0.58%  '

$ws.Range("E16").Value = '  +0.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9998'
$ws.Range("E17").Value = '  +0.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.872.66'
$ws.Range("E18").Value = '  +0.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.71'
$ws.Range("E19").Value = '  +1.49%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006761'
$ws.Range("E20").Value = '  +1.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.966.59'
$ws.Range("E21").Value = '  +0.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.280'
$ws.Range("E22").Value = '  +0.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.639'
$ws.Range("E23").Value = '  -1.95%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.278'
$ws.Range("E24").Value = '  +2.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '138.49'
$ws.Range("E25").Value = '  -1.28%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.516'
$ws.Range("E26").Value = '  +0.33%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.17'
$ws.Range("E27").Value = '  -0.03%  '

$ws.Range("E28").Value = '  -1.41%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '105.62'
$ws.Range("E29").Value = '  +3.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.929'
$ws.Range("E30").Value = '  +5.33%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08284'
$ws.Range("E31").Value = '  -0.56%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.687'
$ws.Range("E32").Value = '  +3.92%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04636'
$ws.Range("E33").Value = '  +2.87%  '

$ws.Range("E34").Value = '  +1.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9885'
$ws.Range("E35").Value = '  +1.65%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6188'
$ws.Range("E36").Value = '  -0.25%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.678'
$ws.Range("E37").Value = '  -0.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01607'
$ws.Range("E38").Value = '  +1.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.943'
$ws.Range("E39").Value = '  +1.76%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9998'
$ws.Range("E40").Value = '  +0.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.20'
$ws.Range("E41").Value = '  -1.79%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.3851'
$ws.Range("E42").Value = '  -0.22%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7398'
$ws.Range("E43").Value = '  +1.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.005'
$ws.Range("E44").Value = '  -0.29%  '

$ws.Range("E45").Value = '  +0.93%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.237'
$ws.Range("E46").Value = '  +0.29%  '

$ws.Range("E47").Value = '  -1.54%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.78'
$ws.Range("E48").Value = '  +2.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.51'
$ws.Range("E49").Value = '  +1.55%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.622'
$ws.Range("E50").Value = '  -0.26%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3409'
$ws.Range("E51").Value = '  +0.17%  '
